# Regeneration of en/fr/es country data books after child program paras updated.
#
# Updates raw input values on the "Programs for children" sheet (rows 2-53,
# columns D-H). The dependent blocks further down the sheet (rows 57-163)
# are driven by shared formulas such as "=D2*0.9" / "=D2*1.05", so Excel's
# automatic recalculation takes care of propagating the new numbers into
# those cells once the inputs below are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Programs for children")

# Row 2
$ws.Range("F2").Value = 0.39473684210526322
$ws.Range("G2").Value = 0.39473684210526322
$ws.Range("H2").Value = 0.39473684210526322

# Row 3
$ws.Range("F3").Value = 0.30769230769230765
$ws.Range("G3").Value = 0.30769230769230765
$ws.Range("H3").Value = 0.30769230769230765

# Row 18
$ws.Range("F18").Value = 0.7

# Row 20
$ws.Range("F20").Value = 0.84

# Row 21
$ws.Range("D21").Value = 0.28260869565217389
$ws.Range("F21").Value = 0

# Row 22
$ws.Range("F22").Value = 0

# Row 23
$ws.Range("D23").Value = 0.28260869565217389
$ws.Range("F23").Value = 0

# Row 24
$ws.Range("F24").Value = 0

# Row 25
$ws.Range("D25").Value = 0.28260869565217389
$ws.Range("F25").Value = 0

# Row 26
$ws.Range("F26").Value = 0

# Row 27
$ws.Range("F27").Value = 1

# Row 28
$ws.Range("F28").Value = 0

# Row 29
$ws.Range("F29").Value = 0

# Row 30
$ws.Range("F30").Value = 1

# Row 31
$ws.Range("F31").Value = 0

# Row 32
$ws.Range("F32").Value = 0

# Row 33
$ws.Range("F33").Value = 1

# Row 34
$ws.Range("F34").Value = 0

# Row 35
$ws.Range("F35").Value = 0

# Row 36
$ws.Range("F36").Value = 1

# Row 37
$ws.Range("F37").Value = 0

# Row 38
$ws.Range("F38").Value = 0

# Row 39
$ws.Range("F39").Value = 1

# Row 40
$ws.Range("F40").Value = 0

# Row 41
$ws.Range("F41").Value = 0

# Row 42
$ws.Range("F42").Value = 0.3

# Row 43
$ws.Range("F43").Value = 0.5

# Row 44
$ws.Range("F44").Value = 0.65

# Row 45
$ws.Range("F45").Value = 0.3

# Row 46
$ws.Range("F46").Value = 0.49

# Row 47
$ws.Range("F47").Value = 0.52

# Row 48
$ws.Range("F48").Value = 0.88

# Row 49
$ws.Range("D49").Value = 0.78409090909090906
$ws.Range("E49").Value = 0.78409090909090906
$ws.Range("F49").Value = 0.78409090909090906
$ws.Range("G49").Value = 0.78409090909090906
$ws.Range("H49").Value = 0.78409090909090906

# Row 50
$ws.Range("D50").Value = 0.88372093023255816
$ws.Range("E50").Value = 0.88372093023255816
$ws.Range("F50").Value = 0.88372093023255816
$ws.Range("G50").Value = 0.88372093023255816
$ws.Range("H50").Value = 0.88372093023255816

# Row 51
$ws.Range("F51").Value = 0.86

# Row 52
$ws.Range("F52").Value = 0

# Row 53
$ws.Range("F53").Value = 0

# --- View/selection bookkeeping -------------------------------------------------
# Move the "Programs for children" view to D2:H53 (also drops the old
# topLeftCell="A100" scroll-freeze since the sheet is no longer scrolled there).
$ws.Range("D2:H53").Select()

# The previously active tab was "Baseline year population inputs"; the
# workbook now opens on "Program dependencies" instead. Activating it last
# makes it the saved active sheet/tab.
$wsDeps = $wb.Worksheets.Item("Program dependencies")
$wsDeps.Activate()
